# Add a reference to RFC 9261 next to the existing draft reference on the
# "Periodic re-authentication" bullet of the Security Requirements slide
# (slide 6, "Content Placeholder 2").
#
# The bullet currently ends with a red run of text reading
# "-post-handshake." We only change that run's text to
# "-post-handshake or RFC 9261." so the run's existing formatting
# (color, font, etc.) is preserved and nothing else on the slide moves.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# Find the content placeholder shape by name (falls back to the known
# positional index if the name ever changes).
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Content Placeholder 2") {
        $shape = $candidate
        break
    }
}
if ($shape -eq $null) {
    $shape = $s.Shapes.Item(2)
}

$tr = $shape.TextFrame.TextRange

$target = "-post-handshake."
$replacement = "-post-handshake or RFC 9261."
$idx = $tr.Text.IndexOf($target)

if ($idx -ge 0) {
    $sub = $tr.Characters($idx + 1, $target.Length)
    $sub.Text = $replacement
}
